# Weekly fruit/vegetable price data refresh.
# Two new daily price rows are inserted at the top of the price block
# (row 815), pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows starting at row 815 (shifts 815.. down by 2).
$ws.Range("A815:R816").Insert()

# --- New row 815 ---
$ws.Range("A815").Value = 10
$ws.Range("B815").Value = "Vega Modelo de Temuco"
$ws.Range("C815").Value = "La Araucanía"
$ws.Range("D815").Value = 45166
$ws.Range("E815").Value = 9
$ws.Range("F815").Value = 100112043
$ws.Range("G815").Value = "Pepino ensalada"
$ws.Range("H815").Value = "Alaska"
$ws.Range("I815").Value = "Primera"
$ws.Range("J815").Value = 80
$ws.Range("K815").Value = 24000
$ws.Range("L815").Value = 24000
$ws.Range("M815").Value = 24000
$ws.Range("N815").Value = "`$/caja 50 unidades"
$ws.Range("O815").Value = "Región de Arica y Parinacota"
$ws.Range("P815").Value = 480
$ws.Range("Q815").Value = 50
$ws.Range("R815").Value = "Hortaliza"

# --- New row 816 ---
$ws.Range("A816").Value = 10
$ws.Range("B816").Value = "Vega Modelo de Temuco"
$ws.Range("C816").Value = "La Araucanía"
$ws.Range("D816").Value = 45166
$ws.Range("E816").Value = 9
$ws.Range("F816").Value = 100112043
$ws.Range("G816").Value = "Pepino ensalada"
$ws.Range("H816").Value = "Sin especificar"
$ws.Range("I816").Value = "Primera"
$ws.Range("J816").Value = 500
$ws.Range("K816").Value = 11000
$ws.Range("L816").Value = 11000
$ws.Range("M816").Value = 11000
$ws.Range("N816").Value = "`$/caja 50 unidades"
$ws.Range("O816").Value = "Región de Arica y Parinacota"
$ws.Range("P816").Value = 220
$ws.Range("Q816").Value = 50
$ws.Range("R816").Value = "Hortaliza"
